# Update cryptocurrency price / volume data to reflect the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) / Volume(1h) (E) refreshes for unchanged coin rows ---
$ws.Range("D2").Value = '44.244.36'
$ws.Range("E2").Value = '  +2.48%  '
$ws.Range("D3").Value = '2.429.47'
$ws.Range("E3").Value = '  +1.93%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '308.42'
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").Value = '100.70'
$ws.Range("E6").Value = '  +3.60%  '
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("D10").Value = '35.41'
$ws.Range("E10").Value = '  +3.40%  '
$ws.Range("D11").Value = '0.0803'
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("E12").Value = '  +2.27%  '
$ws.Range("D13").Value = '18.82'
$ws.Range("E13").Value = '  +1.96%  '
$ws.Range("E14").Value = '  +1.91%  '
$ws.Range("D15").Value = '2.807.77'
$ws.Range("E15").Value = '  +1.98%  '
$ws.Range("D16").Value = '2.425.49'
$ws.Range("E16").Value = '  +1.58%  '
$ws.Range("E17").Value = '  +3.42%  '
$ws.Range("D18").Value = '44.197.88'
$ws.Range("E18").Value = '  +2.38%  '
$ws.Range("D19").Value = '12.35'
$ws.Range("E19").Value = '  +1.36%  '
$ws.Range("E20").Value = '  +1.84%  '
$ws.Range("E21").Value = '  +2.07%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = '241.09'
$ws.Range("E23").Value = '  +2.42%  '
$ws.Range("E24").Value = '  +1.91%  '
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("D27").Value = '25.28'
$ws.Range("E27").Value = '  +1.49%  '
$ws.Range("D28").Value = '2.35'
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("E29").Value = '  +4.93%  '
$ws.Range("D30").Value = '33.23'
$ws.Range("E30").Value = '  +4.92%  '
$ws.Range("E31").Value = '  +13.36%  '
$ws.Range("D32").Value = '18.67'
$ws.Range("E32").Value = '  +8.62%  '
$ws.Range("E33").Value = '  +1.62%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  +3.21%  '
$ws.Range("E36").Value = '  +3.22%  '
$ws.Range("E37").Value = '  +4.80%  '
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("D42").Value = '21.19'
$ws.Range("E42").Value = '  -5.78%  '
$ws.Range("E43").Value = '  +2.66%  '
$ws.Range("D44").Value = '1.957.06'
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("E45").Value = '  +1.83%  '
$ws.Range("D46").Value = '2.89'
$ws.Range("E46").Value = '  +4.94%  '
$ws.Range("D47").Value = '9.44'
$ws.Range("E47").Value = '  +2.51%  '
$ws.Range("D48").Value = '1.67'
$ws.Range("E48").Value = '  +10.34%  '
$ws.Range("D49").Value = '53.51'
$ws.Range("E49").Value = '  +1.29%  '
$ws.Range("D50").Value = '73.79'
$ws.Range("E50").Value = '  +2.55%  '
$ws.Range("E51").Value = '  +0.90%  '

# --- Rows 38/39 swapped rank order (Monero <-> LidoDAOToken) with refreshed data ---
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '2.92'
$ws.Range("E38").Value = '  +4.23%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '128.06'
$ws.Range("E39").Value = '  +24.32%  '
